$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 (pushes old rows 5-13 down to 6-14),
# inheriting formatting from the row above.
$ws.Rows.Item(5).Insert()

# --- Column A: ids 2..13 as text (rows 3..14) ---
$ws.Range("A3").Value  = "2"
$ws.Range("A4").Value  = "3"
$ws.Range("A5").Value  = "4"
$ws.Range("A6").Value  = "5"
$ws.Range("A7").Value  = "6"
$ws.Range("A8").Value  = "7"
$ws.Range("A9").Value  = "8"
$ws.Range("A10").Value = "9"
$ws.Range("A11").Value = "10"
$ws.Range("A12").Value = "11"
$ws.Range("A13").Value = "12"
$ws.Range("A14").Value = "13"

# --- Rename the former "1.1.1 Ist diese Unterlage vollständig?" question ---
$ws.Range("C4").Value = "1.1.1 Einbauanleitung"

# --- Fill in the newly inserted row 5 ---
$ws.Range("B5").Value = "2"
$ws.Range("C5").Value = "1.1.2 Bedienungsanleitung"
$ws.Range("D5").Value = "Nein"

# --- Selection as left by the edit ---
$ws.Range("C6").Select() | Out-Null
